$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3883
$ws.Range("C3").Value = 4091
$ws.Range("C4").Value = 4091
$ws.Range("C5").Value = 4173
$ws.Range("C6").Value = 4332
$ws.Range("C7").Value = 4683
$ws.Range("C8").Value = 4733
$ws.Range("C9").Value = 4733
$ws.Range("C10").Value = 4733
$ws.Range("C11").Value = 4797
$ws.Range("C12").Value = 4797
$ws.Range("C13").Value = 4797
$ws.Range("C14").Value = 4884
$ws.Range("C15").Value = 4953
